$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$changes = @(
    @(2, 374642, 483536102),
    @(3, 297, 359605),
    @(4, 397, 582935),
    @(9, 994, 1553917),
    @(10, 22, 58256),
    @(11, 132994, 205587635),
    @(12, 156, 232177),
    @(13, 70586, 109144113),
    @(17, 4242, 6095466),
    @(22, 9363, 13743227),
    @(24, 88440, 110270639),
    @(25, 59, 104209),
    @(29, 309, 443413),
    @(30, 36066, 54043371),
    @(33, 13208, 19790969),
    @(36, 1668, 2412573),
    @(38, 2455, 3541304),
    @(39, 109639, 137686339),
    @(40, 71, 83868),
    @(41, 91, 127578),
    @(45, 954, 1409096),
    @(47, 49171, 73337646),
    @(49, 10525, 15625112),
    @(51, 1517, 2116660),
    @(54, 3100, 4499618),
    @(55, 79010, 99868132),
    @(57, 59, 97874),
    @(60, 408, 596879),
    @(62, 31812, 47788844),
    @(63, 34, 46616),
    @(65, 12998, 19392153),
    @(67, 1475, 2066523),
    @(71, 1998, 2966860),
    @(73, 23788, 31370244),
    @(77, 8706, 13323743),
    @(79, 6004, 9093621),
    @(80, 583, 840591),
    @(81, 390, 575028),
    @(82, 161606, 202123208),
    @(84, 100, 146097),
    @(86, 482, 711434),
    @(88, 70824, 105892785),
    @(91, 33918, 50743255),
    @(93, 3061, 4487236),
    @(95, 3775, 5551121),
    @(96, 39589, 54500229),
    @(100, 9860, 15112225),
    @(102, 9144, 13844269),
    @(104, 616, 887210),
    @(106, 17261, 33210999),
    @(109, 3971, 8202895),
    @(111, 5677, 11926639),
    @(114, 339, 684578),
    @(116, 164531, 205883880),
    @(120, 1052, 1607748),
    @(122, 59903, 92083873),
    @(123, 112, 172671),
    @(124, 32677, 50380663),
    @(125, 1408, 1968165),
    @(129, 3099, 4609956),
    @(131, 681906, 945223786),
    @(132, 115, 177262),
    @(133, 271, 467901),
    @(136, 1776, 3017512),
    @(137, 41, 89510),
    @(138, 255529, 411152777),
    @(139, 610, 1237405),
    @(140, 25, 62940),
    @(141, 244082, 396185550),
    @(144, 3168, 4586026),
    @(147, 9354, 14121005),
    @(150, 51639, 69561056),
    @(154, 10, 23427),
    @(156, 15771, 23535192),
    @(157, 4307, 6315861),
    @(162, 541, 799761),
    @(163, 20794, 27768776),
    @(167, 8651, 13111102),
    @(169, 6158, 9255703),
    @(172, 351, 522444),
    @(174, 29843, 61901918),
    @(176, 3017, 6284508),
    @(177, 405, 847805),
    @(179, 95, 204537),
    @(180, 191, 421287),
    @(181, 100064, 125400569),
    @(186, 676, 1027881),
    @(188, 38146, 57934875),
    @(190, 15376, 23506230),
    @(192, 1341, 1888554),
    @(194, 2202, 3217879),
    @(196, 271353, 338313219),
    @(198, 198, 284155),
    @(202, 984, 1489026),
    @(204, 96664, 144830901),
    @(207, 38317, 57091266),
    @(210, 5588, 8064460),
    @(213, 6461, 9196944),
    @(216, 305361, 381626217),
    @(223, 687, 1080599),
    @(225, 107878, 165944200),
    @(228, 60349, 92247151),
    @(231, 5055, 7146806),
    @(234, 8155, 11781604),
    @(237, 121893, 152592719),
    @(239, 95, 133401),
    @(240, 15, 20147),
    @(242, 620, 894708),
    @(244, 55437, 82468915),
    @(246, 14798, 22082090),
    @(248, 2002, 2919993),
    @(250, 3279, 4716850),
    @(251, 304243, 391541789),
    @(252, 212, 271833),
    @(253, 285, 404233),
    @(258, 989, 1545978),
    @(260, 111412, 174280746),
    @(263, 79678, 125119576),
    @(265, 2676, 3880032),
    @(268, 6598, 9816671),
)

foreach ($change in $changes) {
    $row = $change[0]
    $newC = $change[1]
    $newD = $change[2]
    $ws.Cells.Item($row, 3).Value = $newC
    $ws.Cells.Item($row, 4).Value = $newD
}
